$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.165.83"
$ws.Range("E2").Value = "  +0.62%  "
$ws.Range("D3").Value = "2.288.93"
$ws.Range("E3").Value = "  +2.05%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "'113.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.50%  "
$ws.Range("D6").Value = "'305.78"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.51%  "
$ws.Range("D7").Value = "'0.633"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.99%  "
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("D10").Value = "'44.78"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.31%  "
$ws.Range("D11").Value = "'0.0929"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.26%  "
$ws.Range("D12").Value = "'55.24"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.38%  "
$ws.Range("D13").Value = "'8.91"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.93%  "
$ws.Range("E14").Value = "  +18.71%  "
$ws.Range("E15").Value = "  -0.38%  "
$ws.Range("D16").Value = "'15.46"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.11%  "
$ws.Range("D17").Value = "2.629.77"
$ws.Range("E17").Value = "  +1.97%  "
$ws.Range("D18").Value = "2.281.23"
$ws.Range("E18").Value = "  +1.72%  "
$ws.Range("D19").Value = "43.071.88"
$ws.Range("E19").Value = "  +0.63%  "
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("D21").Value = "'7.26"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.74%  "
$ws.Range("D22").Value = "'75.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.27%  "
$ws.Range("D23").Value = "'3.58"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +9.22%  "
$ws.Range("D25").Value = "'254.08"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +9.28%  "
$ws.Range("D26").Value = "'9.02"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.17%  "
$ws.Range("D27").Value = "'11.77"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.12%  "
$ws.Range("E28").Value = "  -0.18%  "
$ws.Range("E29").Value = "  +0.17%  "
$ws.Range("D30").Value = "'38.24"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.93%  "
$ws.Range("D31").Value = "'22.30"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.05%  "
$ws.Range("D32").Value = "'175.04"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.25%  "
$ws.Range("E33").Value = "  -3.35%  "
$ws.Range("D34").Value = "'0.0902"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.65%  "
$ws.Range("D35").Value = "'5.69"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.38%  "
$ws.Range("D36").Value = "'5.08"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.48%  "
$ws.Range("E37").Value = "  +0.91%  "
$ws.Range("E38").Value = "  -7.86%  "
$ws.Range("E39").Value = "  +1.41%  "
$ws.Range("E40").Value = "  -0.92%  "
$ws.Range("D41").Value = "'2.56"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.64%  "
$ws.Range("D42").Value = "'72.66"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.56%  "
$ws.Range("E43").Value = "  -0.86%  "
$ws.Range("E44").Value = "  -0.22%  "
$ws.Range("D45").Value = "'12.68"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.80%  "
$ws.Range("D46").Value = "'1.40"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.31%  "
$ws.Range("D47").Value = "'5.67"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.67%  "
$ws.Range("D48").Value = "'107.50"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.96%  "
$ws.Range("D49").Value = "'1.31"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.42%  "
$ws.Range("D50").Value = "'8.83"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.02%  "
$ws.Range("D51").Value = "'74.26"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.95%  "
